$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prêts")

# ---------------------------------------------------------------------
# 1. Re-point the interest-rate inputs: they move from C6/H6/M6 to
#    B6/C6/D6. Copy C6's format (percentage, centered) onto B6 and D6,
#    fill in the three rates, then drop the now-unused H6/M6 cells.
# ---------------------------------------------------------------------
$ws.Range("C6").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("B6").Value = 0.03
$ws.Range("C6").Value = 0.06
$ws.Range("D6").Value = 0.09

$ws.Range("H6").Clear()
$ws.Range("M6").Clear()

# Row 3 formulas now reference the new B6/C6/D6 rate cells instead of
# the old C6/H6/M6 ones (the cell styles are left untouched).
$ws.Range("C3").Formula = "=B3*`$B6+B3"
$ws.Range("D3").Formula = "=C3*`$B6+C3"
$ws.Range("E3").Formula = "=D3*`$B6+D3"
$ws.Range("F3").Formula = "=E3*`$B6+E3"
$ws.Range("G3").Formula = "=F3*`$B6+F3"

$ws.Range("H3").Formula = "=B3+B3*`$C6"
$ws.Range("I3").Formula = "=C3+C3*`$C6"
$ws.Range("J3").Formula = "=D3+D3*`$C6"
$ws.Range("K3").Formula = "=E3+E3*`$C6"
$ws.Range("L3").Formula = "=F3+F3*`$C6"

$ws.Range("M3").Formula = "=B3+B3*`$D6"
$ws.Range("N3").Formula = "=C3+C3*`$D6"
$ws.Range("O3").Formula = "=D3+D3*`$D6"
$ws.Range("P3").Formula = "=E3+E3*`$D6"
$ws.Range("Q3").Formula = "=F3+F3*`$D6"

# ---------------------------------------------------------------------
# 2. Insert the new "standard / risqué" comparison table in rows 8-14.
# ---------------------------------------------------------------------
$ws.Rows("8:14").Insert()

$ws.Range("B8").Value = "SECURITE"
$ws.Range("C8").Value = "TRANQUILITE"
$ws.Range("D8").Value = "DYNAMIQUE"

# B9:D9 reuse the "Milliers" (thousands) format already used elsewhere
# in the workbook (e.g. crédits!M4).
$ws1 = $wb.Worksheets.Item("crédits")
$ws1.Range("M4").Copy()
$ws.Range("B9:D9").PasteSpecial(-4122)
$ws.Range("B9").Value = 10000
$ws.Range("C9").Value = 10000
$ws.Range("D9").Value = 15000

# A10:A14 reuse the centered label format used by row 2's headers.
$ws.Range("C2").Copy()
$ws.Range("A10:A14").PasteSpecial(-4122)
$ws.Range("A10").Value = "Année 1"
$ws.Range("A11").Value = "Année 2"
$ws.Range("A12").Value = "Année 3"
$ws.Range("A13").Value = "Année 4"
$ws.Range("A14").Value = "Année 5"

# B10:D14 reuse the existing euro-with-decimals format (xfId-linked to
# "Milliers") and then get a 2-decimal variant of that format so the
# grown capital displays as e.g. "16 350,00 €".
$ws.Range("L3").Copy()
$ws.Range("B10:D14").PasteSpecial(-4122)
$ws.Range("B10:D14").NumberFormat = '_-* #,##0.00\ [$€-40C]_-;\-* #,##0.00\ [$€-40C]_-;_-* "-"??\ [$€-40C]_-;_-@_-'

$ws.Range("B10").Formula = "=B9*B`$6+B9"
$ws.Range("B11").Formula = "=B10*B`$6+B10"
$ws.Range("B12").Formula = "=B11*B`$6+B11"
$ws.Range("B13").Formula = "=B12*B`$6+B12"
$ws.Range("B14").Formula = "=B13*B`$6+B13"

$ws.Range("C10").Formula = "=C9*C`$6+C9"
$ws.Range("C11").Formula = "=C10*C`$6+C10"
$ws.Range("C12").Formula = "=C11*C`$6+C11"
$ws.Range("C13").Formula = "=C12*C`$6+C12"
$ws.Range("C14").Formula = "=C13*C`$6+C13"

$ws.Range("D10").Formula = "=D9*D`$6+D9"
$ws.Range("D11").Formula = "=D10*D`$6+D10"
$ws.Range("D12").Formula = "=D11*D`$6+D11"
$ws.Range("D13").Formula = "=D12*D`$6+D12"
$ws.Range("D14").Formula = "=D13*D`$6+D13"

# ---------------------------------------------------------------------
# 3. Cosmetic touches: widen column B, clear the frozen/scrolled view,
#    and move the selection to where the user ended up (K25).
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 15.17

$ws.Range("A1").Select()
$ws.Range("K25").Select()
